$d = $word.ActiveDocument
$d.Content.Find.Execute("09/2022 - 10/2023", $true, $false, $false, $false, $false, $true, 1, $false, "09/2022 - 9/2023", 2)
